# Applies the OOXML diff to Filtermatching_GUI.docx:
#  1. Update the italic path text from the 2022 summer-temp folder to the
#     new Filtermatching_GUI/code folder.
#  2. Move the "_GoBack" bookmark from the very end of the document to
#     right after that updated path run (where Word last left the cursor).
#  3/4. Give the two screenshot-containing runs (that previously had no
#     <w:rPr>) an explicit <w:noProof/><w:lang w:eastAsia="nb-NO"/>, matching
#     the rPr already used by the document's other inline drawings.

$d = $word.ActiveDocument

# --- 1. Update the folder path text (contained in a single run) ---------
$oldPath = "F:\Røntgen\Arbeidsmappe\2022\2022 Sommervikar Jostein\programmering"
$newPath = "F:\Røntgen\Arbeidsmappe\2022\2022 Filtermatching_GUI\code"

$findRange = $d.Content
$found = $findRange.Find.Execute($oldPath)
if (-not $found) {
    throw "Could not find the old folder-path text to replace."
}

$pathRange = $d.Range($findRange.Start, $findRange.End)
$pathRange.Text = $newPath

# --- 2. Relocate the "_GoBack" bookmark ----------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$goBackRange = $d.Range($pathRange.End, $pathRange.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# --- 3/4. Add noProof + eastAsia lang to the two screenshot runs --------
foreach ($idx in 9, 10) {
    $shapeRange = $d.InlineShapes.Item($idx).Range
    $shapeRange.NoProofing = $true
    $shapeRange.LanguageIDFarEast = "nb-NO"
}

Write-Output "Edit complete."
